$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously held the "is_locked" label/data-validation template string;
# that template is being dropped, and the column now carries the
# "order_by" comment template that used to live two columns over (F1).
$ws.Range("D1").Value = "<%=comment.order_by%>"

# E1 previously held the "is_enabled" label/data-validation template
# string; that template is being dropped too, and the column now carries
# the "rem" comment template that used to live in G1.
$ws.Range("E1").Value = "<%=comment.rem%>"

# The old F1/G1 cells (order_by/rem) are now redundant since their
# content moved into D1/E1 above, so remove them and shrink the used
# range back down to A1:E1.
$ws.Range("F1:G1").ClearContents() | Out-Null
